$d = $word.ActiveDocument

$d.Paragraphs.Item(1).Range.Text = "⚡️🚀המאמר היומי של מייק 22.06.24:⚡️🚀"
$d.Paragraphs.Item(2).Range.Text = "GLiNER: Generalist Model for Named Entity Recognition using Bidirectional Transformer"
$d.Paragraphs.Item(3).Range.Text = " המאמר הזה הוא שפצור קל של המאמר שסקרנו אתמול 21.06.24. המאמר מציע גישה לאימון והיסק של מודל לזיהוי NER המורכב משלבים הבאים:"
$d.Paragraphs.Item(4).Range.Text = " מעברים כל קטגוריה שברצוננו לזהות דרך טוקנייזר - הקטגוריות מופרדות על ידי טוקן מיוחד הנקרא `"ENT`""
$d.Paragraphs.Item(5).Range.Text = " מעבירים דרך הטוקנייזר את כל הטוקנים של הטקסט. ד״א הטוקנים של הקטגוריות מופרדות מהטוקנים של טקסט על ידי טוקן מיוחד `"SEP`""
$d.Paragraphs.Item(6).Range.Text = "מכניסים את הטוקנים מהשלבים הקודמים לטרנספומר דו-כיווני (encoder) כמו BERT או ROBERTA"
$d.Paragraphs.Item(7).Range.Text = "מעבירים את הייצוגים תלויי הקשר של הקטגוריות דרך FFN דו שכבתי (יש כזה בטרנספורמר) כדי לקבל ייצוג של כל קטגוריה."
$d.Paragraphs.Item(8).Range.Text = "מפעילים את מה שנקרא במאמר הקודם: Structured Span Prediction כלומר כדי לזהות את הקטגוריה של הטוקנים i עד i+n: לוקחים את הייצוג של טוקן ה-i ואת זה של טוקן i+n ומעבירים את השרשור שלהם דרך FFN דו שכבתי (מבנה דומה לסעיף הקודם) וכך מפיקים ייצוגו של ה-span הזה"

# Append the three new trailing paragraphs after paragraph 8 (the URL paragraph)
$null = $d.Paragraphs.Item(8).Range.InsertParagraphAfter()
$d.Paragraphs.Item(9).Range.Text = "כדי לשערך הסתברות ש- span (תת-סדרה של טוקנים רצופים) שייך לקטגוריה j מחשבים סיגמואיד של המכפלה פנימית של ייצוג הקטגוריה j מסעיף 4 עם ייצוג ה-span מהסעיף הקודם."
$null = $d.Paragraphs.Item(9).Range.InsertParagraphAfter()
$d.Paragraphs.Item(10).Range.Text = "מפעילים אלגוריתמיםן גרידים כדי לזהות spans השייכים לכל קטגוריה (המאמר לא מרחיב על כך, צריך להביט בקוד)"
$null = $d.Paragraphs.Item(10).Range.InsertParagraphAfter()
$d.Paragraphs.Item(11).Range.Text = "https://arxiv.org/abs/2311.08526"

Write-Output ("final paragraph count: " + $d.Paragraphs.Count)
